# Update crypto price/volume data per the diff (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether it is a Price-column
# cell that must be forced to Text format so Excel does not reinterpret
# strings like `1.00` or `64.273.61` as numbers/dates.
$updates = @(
    @{ Cell = 'D2'; Value = '64.273.61'; AsText = $true },
    @{ Cell = 'E2'; Value = '  +5.22%  '; AsText = $false },
    @{ Cell = 'D3'; Value = '2.772.75'; AsText = $true },
    @{ Cell = 'E3'; Value = '  +4.58%  '; AsText = $false },
    @{ Cell = 'D4'; Value = '1.00'; AsText = $true },
    @{ Cell = 'E4'; Value = '  -0.14%  '; AsText = $false },
    @{ Cell = 'D5'; Value = '581.04'; AsText = $true },
    @{ Cell = 'E5'; Value = '  +0.62%  '; AsText = $false },
    @{ Cell = 'D6'; Value = '155.09'; AsText = $true },
    @{ Cell = 'E6'; Value = '  +7.05%  '; AsText = $false },
    @{ Cell = 'D7'; Value = '1.00'; AsText = $true },
    @{ Cell = 'E7'; Value = '  +0.10%  '; AsText = $false },
    @{ Cell = 'E8'; Value = '  +1.57%  '; AsText = $false },
    @{ Cell = 'D9'; Value = '2.773.03'; AsText = $true },
    @{ Cell = 'E10'; Value = '  +1.64%  '; AsText = $false },
    @{ Cell = 'D11'; Value = '0.112'; AsText = $true },
    @{ Cell = 'E11'; Value = '  +4.80%  '; AsText = $false },
    @{ Cell = 'D12'; Value = '0.389'; AsText = $true },
    @{ Cell = 'E12'; Value = '  +1.22%  '; AsText = $false },
    @{ Cell = 'E13'; Value = '  +2.88%  '; AsText = $false },
    @{ Cell = 'D14'; Value = '3.261.90'; AsText = $true },
    @{ Cell = 'E14'; Value = '  +4.42%  '; AsText = $false },
    @{ Cell = 'D15'; Value = '26.58'; AsText = $true },
    @{ Cell = 'D16'; Value = '64.199.68'; AsText = $true },
    @{ Cell = 'E16'; Value = '  +5.11%  '; AsText = $false },
    @{ Cell = 'E17'; Value = '  +5.83%  '; AsText = $false },
    @{ Cell = 'D18'; Value = '2.771.60'; AsText = $true },
    @{ Cell = 'E18'; Value = '  +3.85%  '; AsText = $false },
    @{ Cell = 'D19'; Value = '12.02'; AsText = $true },
    @{ Cell = 'E19'; Value = '  +3.05%  '; AsText = $false },
    @{ Cell = 'E20'; Value = '  +2.65%  '; AsText = $false },
    @{ Cell = 'D21'; Value = '361.69'; AsText = $true },
    @{ Cell = 'E21'; Value = '  +2.86%  '; AsText = $false },
    @{ Cell = 'D22'; Value = '7.04'; AsText = $true },
    @{ Cell = 'E22'; Value = '  +1.34%  '; AsText = $false },
    @{ Cell = 'E23'; Value = '  +0.23%  '; AsText = $false },
    @{ Cell = 'D24'; Value = '0.533'; AsText = $true },
    @{ Cell = 'E24'; Value = '  +0.86%  '; AsText = $false },
    @{ Cell = 'D25'; Value = '66.59'; AsText = $true },
    @{ Cell = 'E25'; Value = '  +4.03%  '; AsText = $false },
    @{ Cell = 'E26'; Value = '  +5.90%  '; AsText = $false },
    @{ Cell = 'D27'; Value = '8.51'; AsText = $true },
    @{ Cell = 'E27'; Value = '  +4.47%  '; AsText = $false },
    @{ Cell = 'E28'; Value = '  +0.20%  '; AsText = $false },
    @{ Cell = 'D29'; Value = '0.0₃0905'; AsText = $true },
    @{ Cell = 'E29'; Value = '  +10.79%  '; AsText = $false },
    @{ Cell = 'D30'; Value = '2.01'; AsText = $true },
    @{ Cell = 'E30'; Value = '  +1.96%  '; AsText = $false },
    @{ Cell = 'E31'; Value = '  +2.78%  '; AsText = $false },
    @{ Cell = 'D32'; Value = '1.30'; AsText = $true },
    @{ Cell = 'E32'; Value = '  +20.46%  '; AsText = $false },
    @{ Cell = 'D33'; Value = '171.61'; AsText = $true },
    @{ Cell = 'E33'; Value = '  +3.69%  '; AsText = $false },
    @{ Cell = 'D34'; Value = '0.999'; AsText = $true },
    @{ Cell = 'E34'; Value = '  +0.03%  '; AsText = $false },
    @{ Cell = 'E35'; Value = '  +2.85%  '; AsText = $false },
    @{ Cell = 'E36'; Value = '  +7.39%  '; AsText = $false },
    @{ Cell = 'D37'; Value = '1.43'; AsText = $true },
    @{ Cell = 'E37'; Value = '  +8.54%  '; AsText = $false },
    @{ Cell = 'E38'; Value = '  +9.25%  '; AsText = $false },
    @{ Cell = 'D39'; Value = '1.01'; AsText = $true },
    @{ Cell = 'E39'; Value = '  +13.69%  '; AsText = $false },
    @{ Cell = 'D40'; Value = '347.07'; AsText = $true },
    @{ Cell = 'E40'; Value = '  +2.91%  '; AsText = $false },
    @{ Cell = 'D41'; Value = '4.24'; AsText = $true },
    @{ Cell = 'E41'; Value = '  +5.02%  '; AsText = $false },
    @{ Cell = 'D42'; Value = '39.25'; AsText = $true },
    @{ Cell = 'E42'; Value = '  +1.62%  '; AsText = $false },
    @{ Cell = 'D43'; Value = '5.82'; AsText = $true },
    @{ Cell = 'E43'; Value = '  +10.98%  '; AsText = $false },
    @{ Cell = 'D44'; Value = '21.89'; AsText = $true },
    @{ Cell = 'E44'; Value = '  +7.13%  '; AsText = $false },
    @{ Cell = 'D45'; Value = '21.92'; AsText = $true },
    @{ Cell = 'E45'; Value = '  +6.22%  '; AsText = $false },
    @{ Cell = 'D46'; Value = '0.651'; AsText = $true },
    @{ Cell = 'E46'; Value = '  +5.67%  '; AsText = $false },
    @{ Cell = 'D47'; Value = '0.0591'; AsText = $true },
    @{ Cell = 'E47'; Value = '  +5.21%  '; AsText = $false },
    @{ Cell = 'D48'; Value = '137.84'; AsText = $true },
    @{ Cell = 'E48'; Value = '  +2.60%  '; AsText = $false },
    @{ Cell = 'D49'; Value = '0.0257'; AsText = $true },
    @{ Cell = 'E49'; Value = '  +3.04%  '; AsText = $false },
    @{ Cell = 'E50'; Value = '  +1.20%  '; AsText = $false },
    @{ Cell = 'D51'; Value = '1.00'; AsText = $true },
    @{ Cell = 'E51'; Value = '  +0.23%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

